$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the contents of columns P (16) and Q (17) ---
# The two columns ("assColl" / occurrence numbers and "assCollTaxa" / taxa
# names) were transposed: every row that had data in P or Q gets those
# two cell values swapped.
$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $pCell = $ws.Cells.Item($r, 16)
    $qCell = $ws.Cells.Item($r, 17)
    $pVal = $pCell.Value2
    $qVal = $qCell.Value2
    if (($pVal -ne $null) -or ($qVal -ne $null)) {
        $pCell.Value2 = $qVal
        $qCell.Value2 = $pVal
    }
}

# --- Swap the column widths for P and Q to match the moved content ---
$pWidth = $ws.Columns.Item(16).ColumnWidth
$qWidth = $ws.Columns.Item(17).ColumnWidth
$ws.Columns.Item(16).ColumnWidth = $qWidth
$ws.Columns.Item(17).ColumnWidth = $pWidth

# --- Swap the header comments attached to P1 / Q1 so they describe the
#     column that now actually lives there ---
$pComment = $ws.Range("P1").Comment
$qComment = $ws.Range("Q1").Comment
$pText = $pComment.Text()
$qText = $qComment.Text()
$pComment.Text($qText)
$qComment.Text($pText)

# --- Update the active selection to reflect where editing left off ---
$ws.Range("Q209").Select()
